$d = $word.ActiveDocument

# Locate the insertion point: right after "Ministerio de Salud, " and
# right before "aprueba" in the first CONSIDERANDO paragraph.
$text = $d.Content.Text
$anchor = "Ministerio de Salud, aprueba"
$idx = $text.IndexOf($anchor)
$insertPos = $idx + 21  # length of "Ministerio de Salud, "

$newText = "sus modificaciones o aquella que la reemplace, se "

# Insert the new text as its own run right before "aprueba".
$r = $d.Range($insertPos, $insertPos)
$r.InsertBefore($newText)

# The engine coalesces adjacent runs that share identical formatting.
# Nudge formatting (toggle Bold on/off, a no-visual-effect round trip)
# on the three affected spans so each keeps its own run, matching the
# surgical single-run insertion described by the change:
#   ","  |  " "  |  <new run>  |  "aprueba"  |  " el \""

# the lone space run right before the insertion point
$rSpace = $d.Range($insertPos - 1, $insertPos)
$rSpace.Bold = 1
$rSpace.Bold = 0

# the newly-inserted text
$rNew = $d.Range($insertPos, $insertPos + $newText.Length)
$rNew.Bold = 1
$rNew.Bold = 0

# "aprueba" immediately following the insertion
$apruebaStart = $insertPos + $newText.Length
$apruebaEnd = $apruebaStart + 7
$rAprueba = $d.Range($apruebaStart, $apruebaEnd)
$rAprueba.Bold = 1
$rAprueba.Bold = 0
